# Adds price-monitor log rows 162-181 (2024-09-30 entries) to Sheet1,
# matching the appended rows in the target OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 162; A = '2024-09-30 17:22:22'; B = 'monitor_price'; C = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'; D = '$69.99'; E = '2024-09-30'; F = '17:22:22' },
    @{ Row = 163; A = '2024-09-30 17:25:28'; B = 'monitor_price'; C = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'; D = '$69.99'; E = '2024-09-30'; F = '17:25:28' },
    @{ Row = 164; A = '2024-09-30 17:25:56'; B = 'monitor_price'; C = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'; D = '$69.99'; E = '2024-09-30'; F = '17:25:56' },
    @{ Row = 165; A = '2024-09-30 17:32:24'; B = 'monitor_price'; C = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'; D = '$69.99'; E = '2024-09-30'; F = '17:32:24' },
    @{ Row = 166; A = '2024-09-30 17:32:56'; B = 'monitor_price'; C = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'; D = '$69.99'; E = '2024-09-30'; F = '17:32:56' },
    @{ Row = 167; A = '2024-09-30 18:06:42'; B = 'monitor_price'; C = 'https://www.ebay.com/itm/314411766963?_trkparms=amclksrc%3DITM%26aid%3D777008%26algo%3DPERSONAL.TOPIC%26ao%3D1%26asc%3D20240603121456%26meid%3Da07931f944bc4a5b95376fe64d0ab035%26pid%3D102177%26rk%3D1%26rkt%3D1%26itm%3D314411766963%26pmt%3D1%26noa%3D1%26pg%3D4375194%26algv%3DNoSignalMostWatched%26brand%3DSimpliSafe&_trksid=p4375194.c102177.m166540&_trkparms=parentrq%3A71497a9c1910a8cd54f819a0ffff582e%7Cpageci%3A59d1354a-5f2b-11ef-9c4d-f2c982e61003%7Ciid%3A1%7Cvlpname%3Avlp_homepage'; D = 'US $179.99'; E = '2024-09-30'; F = '18:06:42' },
    @{ Row = 168; A = '2024-09-30 18:07:08'; B = 'monitor_price'; C = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'; D = '$69.99'; E = '2024-09-30'; F = '18:07:08' },
    @{ Row = 169; A = '2024-09-30 18:07:32'; B = 'monitor_price'; C = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'; D = '$69.99'; E = '2024-09-30'; F = '18:07:32' },
    @{ Row = 170; A = '2024-09-30 18:07:56'; B = 'monitor_price'; C = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'; D = '$69.99'; E = '2024-09-30'; F = '18:07:56' },
    @{ Row = 171; A = '2024-09-30 18:08:19'; B = 'monitor_price'; C = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'; D = '$69.99'; E = '2024-09-30'; F = '18:08:19' },
    @{ Row = 172; A = '2024-09-30 18:09:24'; B = 'monitor_price'; C = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'; D = '$69.99'; E = '2024-09-30'; F = '18:09:24' },
    @{ Row = 173; A = '2024-09-30 21:51:31'; B = 'monitor_price'; C = 'https://example.com/product'; D = '$199.99'; E = '2024-09-30'; F = '21:51:31' },
    @{ Row = 174; A = '2024-09-30 21:51:31'; B = 'monitor_price'; C = 'invalid_url'; D = 'Error fetching price: Invalid URL'; E = '2024-09-30'; F = '21:51:31' },
    @{ Row = 175; A = '2024-09-30 21:51:34'; B = 'monitor_price'; C = 'https://example.com/product'; D = '100 USD'; E = '2024-09-30'; F = '21:51:34' },
    @{ Row = 176; A = '2024-09-30 21:54:07'; B = 'monitor_price'; C = 'https://example.com/product'; D = '$199.99'; E = '2024-09-30'; F = '21:54:07' },
    @{ Row = 177; A = '2024-09-30 21:54:07'; B = 'monitor_price'; C = 'invalid_url'; D = 'Error fetching price: Invalid URL'; E = '2024-09-30'; F = '21:54:07' },
    @{ Row = 178; A = '2024-09-30 21:54:10'; B = 'monitor_price'; C = 'https://example.com/product'; D = '100 USD'; E = '2024-09-30'; F = '21:54:10' },
    @{ Row = 179; A = '2024-09-30 22:07:53'; B = 'monitor_price'; C = 'https://example.com/product'; D = '$199.99'; E = '2024-09-30'; F = '22:07:53' },
    @{ Row = 180; A = '2024-09-30 22:07:53'; B = 'monitor_price'; C = 'invalid_url'; D = 'Error fetching price: Invalid URL'; E = '2024-09-30'; F = '22:07:53' },
    @{ Row = 181; A = '2024-09-30 22:07:55'; B = 'monitor_price'; C = 'https://example.com/product'; D = '100 USD'; E = '2024-09-30'; F = '22:07:55' }
)

foreach ($r in $newRows) {
    foreach ($col in @("A", "B", "C", "D", "E", "F")) {
        $addr = "$col$($r.Row)"
        # Leading apostrophe forces Excel to store the value as literal text
        # (matching t="inlineStr"/text cells in the source log) instead of
        # auto-detecting a number, currency, or date/time value.
        $ws.Range($addr).Value = "'" + $r.$col
        $ws.Range($addr).Style = "Normal"
    }
}
